$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the product group for rows 33-34 from "PRoduct" to "Crispril"
$ws.Range("A33").Value = "Crispril"
$ws.Range("A34").Value = "Crispril"

# Give these two cells a (bottom) border - this is what introduces the
# new cellXfs entries with applyBorder="1" in styles.xml
$ws.Range("A33").Borders.Item(9).LineStyle = 1
$ws.Range("A34").Borders.Item(9).LineStyle = 1

# Update the view: scroll down and move the active selection
$ws.Range("A35").Select()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
